$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.703.03"
$ws.Range("E2").Value = "  +2.88%  "

$ws.Range("D3").Value = "3.472.50"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.29"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.17"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +1.63%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.787"
$ws.Range("E9").Value = "  +8.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +20.04%  "

$ws.Range("B11").Value = "ShibaInu"
$ws.Range("C11").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000273"
$ws.Range("E11").Value = "  +28.45%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.47"
$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.85"
$ws.Range("E13").Value = "  +7.61%  "

$ws.Range("D14").Value = "4.018.41"
$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.35"
$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("D17").Value = "3.455.14"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.09"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.41"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").Value = "63.495.18"
$ws.Range("E20").Value = "  +2.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "462.15"
$ws.Range("E21").Value = "  -4.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.27"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.62"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.23"
$ws.Range("E25").Value = "  +11.25%  "

$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.17"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.77"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.43"
$ws.Range("E29").Value = "  +3.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.53"
$ws.Range("E30").Value = "  -3.43%  "

$ws.Range("E31").Value = "  -0.52%  "

$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("E33").Value = "  -0.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.98"
$ws.Range("E34").Value = "  -4.46%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.71"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0489"
$ws.Range("E37").Value = "  -0.79%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.143"
$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +4.60%  "

$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("D42").Value = "0.0₃0655"
$ws.Range("E42").Value = "  +60.17%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.50"
$ws.Range("E43").Value = "  +6.69%  "

$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.33"
$ws.Range("E44").Value = "  -2.32%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.317"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.01"
$ws.Range("E46").Value = "  -5.07%  "

$ws.Range("E47").Value = "  -6.01%  "

$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.75"
$ws.Range("E49").Value = "  -4.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.91"
$ws.Range("E50").Value = "  -3.27%  "

$ws.Range("E51").Value = "  -7.49%  "
